# Update with restock suggestion
# Applies forecast-refresh edits to "Forecast Comparison" and "Summary" sheets.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 1: "Forecast Comparison"
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Week_Start_Date (column B) for each data row (weeks W1..W16, rows 2..17)
$weekDates = @{
    2  = "2023-09-03"
    3  = "2023-09-10"
    4  = "2023-09-17"
    5  = "2023-09-24"
    6  = "2023-10-01"
    7  = "2023-10-08"
    8  = "2023-10-15"
    9  = "2023-10-22"
    10 = "2023-10-29"
    11 = "2023-11-05"
    12 = "2023-11-12"
    13 = "2023-11-19"
    14 = "2023-11-26"
    15 = "2023-12-03"
    16 = "2023-12-10"
    17 = "2023-12-17"
}

# Updated Seasonality Index (column P) for each row
$seasonality = @{
    2  = 1.18
    3  = 0.95
    4  = 1.03
    5  = 1.08
    6  = 0.89
    7  = 1.14
    8  = 0.92
    9  = 1.14
    10 = 1.06
    11 = 1.07
    12 = 0.81
    13 = 0.96
    14 = 1.03
    15 = 0.97
    16 = 1.08
    17 = 1.1
}

# Column Q header becomes "Lifecycle Stage" (previously "Sales Volume Rank")
$ws1.Range("Q1").Value = "'Lifecycle Stage"
$ws1.Range("Q1").Style = "Normal"

for ($row = 2; $row -le 17; $row++) {
    # Week_Start_Date: was blank, now populated. Force text so the
    # date-like string is not auto-converted to a date serial.
    $cellB = $ws1.Cells.Item($row, 2)
    $cellB.Value = "'" + $weekDates[$row]
    $cellB.Style = "Normal"

    # Seasonality Index: plain numeric update.
    $ws1.Cells.Item($row, 16).Value = $seasonality[$row]

    # Column Q now carries the Lifecycle Stage text ("Decline") instead of
    # the old numeric Sales Volume Rank value.
    $cellQ = $ws1.Cells.Item($row, 17)
    $cellQ.Value = "'Decline"
    $cellQ.Style = "Normal"
}

# Rows for W7 (row 8) and W16 (row 17) also get revised restock fields:
#   Inventory Coverage -> blank, Stockout Risk -> Low, Reorder Urgency -> Normal
foreach ($row in @(8, 17)) {
    $ws1.Cells.Item($row, 12).Value = $null

    $cellM = $ws1.Cells.Item($row, 13)
    $cellM.Value = "'Low"
    $cellM.Style = "Normal"

    $cellN = $ws1.Cells.Item($row, 14)
    $cellN.Value = "'Normal"
    $cellN.Style = "Normal"
}

# Column R ("Lifecycle Stage", now superseded by column Q) is removed entirely.
$ws1.Columns("R").Delete()

# -----------------------------------------------------------------
# Sheet 2: "Summary"
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$summaryValues = @{
    2  = "2022-12-25 to 2024-12-29"
    4  = "19"
    5  = "6"
    6  = "2"
    7  = "7"
    8  = "117 units"
    9  = "12"
    10 = "7"
    11 = "4"
    12 = "1"
    13 = "N/A"
    14 = "0"
    15 = "N/A"
}

foreach ($row in $summaryValues.Keys) {
    $cell = $ws2.Cells.Item($row, 2)
    $cell.Value = "'" + $summaryValues[$row]
    $cell.Style = "Normal"
}
